# The "category" and "group" columns (D/E = name pair, F/G = code pair)
# were mislabeled/swapped in the source data. This fixes it by swapping
# the contents of column D <-> column E, and column F <-> column G,
# across the full used range (including the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$firstRow = $used.Row

$rangeDE = $ws.Range("D$($firstRow):E$($lastRow)")
$valsDE = $rangeDE.Value2

$rangeFG = $ws.Range("F$($firstRow):G$($lastRow)")
$valsFG = $rangeFG.Value2

$rowCount = $valsDE.GetLength(0)

for ($i = 1; $i -le $rowCount; $i++) {
    $d = $valsDE[$i,1]
    $e = $valsDE[$i,2]
    $valsDE[$i,1] = $e
    $valsDE[$i,2] = $d

    $f = $valsFG[$i,1]
    $g = $valsFG[$i,2]
    $valsFG[$i,1] = $g
    $valsFG[$i,2] = $f
}

$rangeDE.Value2 = $valsDE
$rangeFG.Value2 = $valsFG
